$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures (row 11, row 13) ---
# "VALOR MORA" total (sum of column F across the data rows)
$ws.Range("E11").Value = 399826
# "Cant. Trabajadores" (distinct workers) and "Cant. Periodos" (data rows)
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 8

# --- Rewrite the worker/period detail table (rows 16-23) ---
# Row 16: DIANA CECILIA YEPES AREVALO (unchanged)
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "45760855"
$ws.Range("D16").Value = "DIANA CECILIA YEPES AREVALO"
$ws.Range("E16").Value = "1811"
$ws.Range("F16").Value = 38000
$ws.Range("G16").Value = 950000

# Rows 17-20: LAURA TATIANA CASTAÑO PARDO, periods 1902-1905
$lauraPeriods = @("1902", "1903", "1904", "1905")
for ($i = 0; $i -lt $lauraPeriods.Length; $i++) {
    $r = 17 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "1143396138"
    $ws.Range("D$r").Value = "LAURA TATIANA CASTAÑO PARDO"
    $ws.Range("E$r").Value = $lauraPeriods[$i]
    $ws.Range("F$r").Value = 33200
    $ws.Range("G$r").Value = 830000
}

# Rows 21-23: MARY CRUZ FONTALVO ORDOÑEZ, periods 2207-2209
$maryPeriods = @("2207", "2208", "2209")
for ($i = 0; $i -lt $maryPeriods.Length; $i++) {
    $r = 21 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "1052997984"
    $ws.Range("D$r").Value = "MARY CRUZ FONTALVO ORDOÑEZ"
    $ws.Range("E$r").Value = $maryPeriods[$i]
    $ws.Range("F$r").Value = 76342
    $ws.Range("G$r").Value = 1908526
}

# Row 24 carried the special "bottom of table" border formatting (it was the last
# data row before the edit). Since the table now ends at row 23, copy that border
# formatting onto row 23 before row 24 is removed.
$ws.Range("B24:J24").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 24 is removed entirely, shifting the signature block (rows 29-30) up to 28-29.
$ws.Rows("24").Delete()
